$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.831.73"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").Value = "1.643.13"
$ws.Range("E3").Value = "  -0.12%  "
$ws.Range("E4").Value = "  -0.35%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.59"
$ws.Range("E5").Value = "  +0.94%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.501"
$ws.Range("E6").Value = "  +0.34%  "
$ws.Range("E7").Value = "  -0.32%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("E9").Value = "  -0.92%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.26"
$ws.Range("E10").Value = "  +0.42%  "
$ws.Range("E11").Value = "  +0.84%  "
$ws.Range("D12").Value = "1.872.27"
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("D13").Value = "1.644.29"
$ws.Range("E13").Value = "  -0.45%  "
$ws.Range("E14").Value = "  -0.21%  "
$ws.Range("E15").Value = "  -0.38%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.33"
$ws.Range("E16").Value = "  +1.46%  "
$ws.Range("D17").Value = "26.840.69"
$ws.Range("E17").Value = "  +0.15%  "
$ws.Range("D18").Value = "0.0₃0735"
$ws.Range("E18").Value = "  -0.43%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "215.92"
$ws.Range("E19").Value = "  +1.04%  "
$ws.Range("E20").Value = "  -0.40%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.67"
$ws.Range("E21").Value = "  +6.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.37"
$ws.Range("E22").Value = "  +0.25%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.37"
$ws.Range("E23").Value = "  -0.61%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.21"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.55"
$ws.Range("E25").Value = "  +1.75%  "
$ws.Range("E26").Value = "  -0.48%  "
$ws.Range("E27").Value = "  +0.35%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.15"
$ws.Range("E28").Value = "  +0.97%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.76"
$ws.Range("E29").Value = "  +0.56%  "
$ws.Range("E30").Value = "  -0.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.20"
$ws.Range("E31").Value = "  +1.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.38"
$ws.Range("E32").Value = "  +2.13%  "
$ws.Range("E33").Value = "  -0.48%  "
$ws.Range("E34").Value = "  +1.06%  "
$ws.Range("D35").Value = "1.270.02"
$ws.Range("E35").Value = "  -1.29%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.45"
$ws.Range("E36").Value = "  +0.39%  "
$ws.Range("E37").Value = "  +0.91%  "
$ws.Range("E38").Value = "  -0.92%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.817"
$ws.Range("E39").Value = "  -1.05%  "
$ws.Range("E40").Value = "  -0.32%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.804"
$ws.Range("E41").Value = "  -0.59%  "
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("D43").Value = "1.783.23"
$ws.Range("E43").Value = "  -0.66%  "
$ws.Range("E44").Value = "  -4.53%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.76"
$ws.Range("E45").Value = "  +1.53%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "61.21"
$ws.Range("E46").Value = "  +0.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.61"
$ws.Range("E47").Value = "  +0.34%  "
$ws.Range("E48").Value = "  -1.63%  "
$ws.Range("E49").Value = "  -0.44%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.54"
$ws.Range("E50").Value = "  -1.64%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0967"
$ws.Range("E51").Value = "  -1.06%  "
